# Iraq_FX.xlsx update: correct row 215 (Aug-2023) close/high values and
# append three new monthly rows (216-218) for Sep, Oct, Nov 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 215 -------------------------------------------------
$ws.Range("D215").Value = 1309
$ws.Range("F215").Value = 1308.62

# --- Row 216 (2023-09-01) ---------------------------------------------------
# Copy the date format (style) from the row above so the new date cell keeps
# the same number format (YYYY-MM-DD HH:MM:SS) as the rest of column A.
$ws.Range("A215").Copy()
$ws.Range("A216:A218").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A216").Value = 45170.33333333334
$ws.Range("B216").Value = "FX_IDC:USDIQD"
$ws.Range("C216").Value = 1308.62
$ws.Range("D216").Value = 1308.65
$ws.Range("E216").Value = 1307
$ws.Range("F216").Value = 1308.62
$ws.Range("G216").Value = 0

# --- Row 217 (2023-10-02) ---------------------------------------------------
$ws.Range("A217").Value = 45201.375
$ws.Range("B217").Value = "FX_IDC:USDIQD"
$ws.Range("C217").Value = 1308.62
$ws.Range("D217").Value = 1310
$ws.Range("E217").Value = 1306
$ws.Range("F217").Value = 1309
$ws.Range("G217").Value = 0

# --- Row 218 (2023-11-01) ---------------------------------------------------
$ws.Range("A218").Value = 45231.375
$ws.Range("B218").Value = "FX_IDC:USDIQD"
$ws.Range("C218").Value = 1309
$ws.Range("D218").Value = 1309
$ws.Range("E218").Value = 1308
$ws.Range("F218").Value = 1309
$ws.Range("G218").Value = 0
